$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.333.18'
$ws.Range("E2").Value = '  -2.33%  '

$ws.Range("D3").Value = '1.793.04'
$ws.Range("E3").Value = '  -2.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.55%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  -0.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.74'
$ws.Range("E6").Value = '  -1.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4507'
$ws.Range("E7").Value = '  -1.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3592'
$ws.Range("E8").Value = '  -3.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.76'
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07083'
$ws.Range("E10").Value = '  -1.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8838'
$ws.Range("E11").Value = '  +0.68%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07748'
$ws.Range("E12").Value = '  -0.64%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.44'
$ws.Range("E13").Value = '  -1.12%  '

$ws.Range("D14").Value = '1.785.25'
$ws.Range("E14").Value = '  -2.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.276'
$ws.Range("E15").Value = '  -1.22%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.323'
$ws.Range("E16").Value = '  -1.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.84'
$ws.Range("E17").Value = '  -2.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  -0.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008506'
$ws.Range("E19").Value = '  -2.46%  '

$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.24'
$ws.Range("E21").Value = '  -1.88%  '

$ws.Range("D22").Value = '26.359.68'
$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.970'
$ws.Range("E23").Value = '  -0.86%  '

$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.033.93'
$ws.Range("E24").Value = '  -0.93%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.52'
$ws.Range("E25").Value = '  +0.83%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.974'
$ws.Range("E26").Value = '  -2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.41'
$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("E28").Value = '  -2.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.017'
$ws.Range("E29").Value = '  +2.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '111.89'
$ws.Range("E30").Value = '  -1.94%  '

$ws.Range("E31").Value = '  -1.07%  '

$ws.Range("E32").Value = '  -1.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.066'
$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.742'
$ws.Range("E34").Value = '  +6.99%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.443'
$ws.Range("E35").Value = '  -0.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7229'
$ws.Range("E36").Value = '  -3.85%  '

$ws.Range("E37").Value = '  -2.88%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.003'
$ws.Range("E38").Value = '  -0.11%  '

$ws.Range("E39").Value = '  -2.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01928'
$ws.Range("E40").Value = '  -1.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05086'
$ws.Range("E41").Value = '  -1.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.854'
$ws.Range("E42").Value = '  -1.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5069'
$ws.Range("E43").Value = '  +1.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.854'
$ws.Range("E44").Value = '  -1.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1516'
$ws.Range("E45").Value = '  -5.46%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.995'
$ws.Range("E46").Value = '  -3.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("E47").Value = '  -0.40%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4629'
$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.14'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.810'
$ws.Range("E50").Value = '  -3.19%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.576'
$ws.Range("E51").Value = '  -2.33%  '
